$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 507.1111
$ws.Cells.Item(28, 9).Value = 463.4091
$ws.Cells.Item(28, 11).Value = 463.4091
$ws.Cells.Item(28, 13).Value = 21.59089999999998
$ws.Cells.Item(33, 8).Value = 330.77777
$ws.Cells.Item(33, 9).Value = 309.625
$ws.Cells.Item(33, 10).Value = 500
$ws.Cells.Item(33, 11).Value = 309.625
$ws.Cells.Item(33, 12).Value = 500
$ws.Cells.Item(33, 13).Value = -80.625
$ws.Cells.Item(33, 14).Value = -958
$ws.Cells.Item(76, 8).Value = 3209.4
$ws.Cells.Item(76, 10).Value = 3700.8
$ws.Cells.Item(76, 12).Value = 3700.8
$ws.Cells.Item(76, 14).Value = -4330.8
$ws.Cells.Item(79, 8).Value = 3209.4
$ws.Cells.Item(79, 10).Value = 3700.8
$ws.Cells.Item(79, 12).Value = 3700.8
$ws.Cells.Item(79, 14).Value = -5884.8
$ws.Cells.Item(92, 8).Value = 625.2308
$ws.Cells.Item(92, 9).Value = 592.8
$ws.Cells.Item(92, 10).Value = 733.3333
$ws.Cells.Item(92, 11).Value = 592.8
$ws.Cells.Item(92, 12).Value = 733.3333
$ws.Cells.Item(92, 13).Value = 655.2
$ws.Cells.Item(92, 14).Value = -3229.3333
$ws.Cells.Item(95, 8).Value = 79960.664
$ws.Cells.Item(95, 10).Value = 79960.664
$ws.Cells.Item(95, 12).Value = 79960.664
$ws.Cells.Item(95, 14).Value = -85452.664
$ws.Cells.Item(137, 8).Value = 4555.8296
$ws.Cells.Item(137, 9).Value = 1286.1923
$ws.Cells.Item(137, 10).Value = 8603.951999999999
$ws.Cells.Item(137, 11).Value = 3858.5769
$ws.Cells.Item(137, 12).Value = 25811.856
$ws.Cells.Item(137, 13).Value = -1308.5769
$ws.Cells.Item(137, 14).Value = -30911.856
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(21, 8).Value = 1063
$ws.Cells.Item(21, 10).Value = 500
$ws.Cells.Item(21, 12).Value = 500
$ws.Cells.Item(21, 14).Value = -1248
$ws.Cells.Item(45, 8).Value = 4337.5
$ws.Cells.Item(45, 10).Value = 5133.3335
$ws.Cells.Item(45, 12).Value = 5133.3335
$ws.Cells.Item(45, 14).Value = -5887.3335
$ws.Cells.Item(55, 8).Value = 33026
$ws.Cells.Item(55, 10).Value = 33026
$ws.Cells.Item(55, 12).Value = 33026
$ws.Cells.Item(55, 14).Value = -33656
$ws.Cells.Item(63, 10).Value = 2000
$ws.Cells.Item(63, 12).Value = 2000
$ws.Cells.Item(63, 14).Value = -3372
$ws.Cells.Item(66, 10).Value = 2000
$ws.Cells.Item(66, 12).Value = 10000
$ws.Cells.Item(66, 14).Value = -16864
$ws.Cells.Item(74, 8).Value = 4378.387
$ws.Cells.Item(74, 9).Value = 897.3674
$ws.Cells.Item(74, 10).Value = 17499.154
$ws.Cells.Item(74, 11).Value = 897.3674
$ws.Cells.Item(74, 12).Value = 17499.154
$ws.Cells.Item(74, 13).Value = -23.36739999999998
$ws.Cells.Item(74, 14).Value = -19247.154
$ws.Cells.Item(77, 8).Value = 4378.387
$ws.Cells.Item(77, 9).Value = 897.3674
$ws.Cells.Item(77, 10).Value = 17499.154
$ws.Cells.Item(77, 11).Value = 4486.837
$ws.Cells.Item(77, 12).Value = 87495.76999999999
$ws.Cells.Item(77, 13).Value = -118.8369999999995
$ws.Cells.Item(77, 14).Value = -96231.76999999999
$ws.Cells.Item(88, 8).Value = 2821.75
$ws.Cells.Item(88, 10).Value = 3007.8572
$ws.Cells.Item(88, 12).Value = 3007.8572
$ws.Cells.Item(88, 14).Value = -3819.8572
$ws.Cells.Item(91, 8).Value = 2821.75
$ws.Cells.Item(91, 10).Value = 3007.8572
$ws.Cells.Item(91, 12).Value = 3007.8572
$ws.Cells.Item(91, 14).Value = -5815.8572
$ws.Cells.Item(104, 8).Value = 29245
$ws.Cells.Item(104, 10).Value = 29245
$ws.Cells.Item(104, 12).Value = 29245
$ws.Cells.Item(104, 14).Value = -36233
$ws.Cells.Item(122, 8).Value = 1719.8
$ws.Cells.Item(122, 9).Value = 1349.8
$ws.Cells.Item(122, 11).Value = 4049.4
$ws.Cells.Item(122, 13).Value = -1599.4
$ws.Cells.Item(132, 8).Value = 642746.3
$ws.Cells.Item(132, 9).Value = 677297.7
$ws.Cells.Item(132, 11).Value = 2031893.1
$ws.Cells.Item(132, 13).Value = -2029363.1
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 32269.559
$ws.Cells.Item(20, 9).Value = 49210.137
$ws.Cells.Item(20, 11).Value = 49210.137
$ws.Cells.Item(20, 13).Value = -48963.137
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 13).ClearContents()
$ws.Cells.Item(105, 8).Value = 2572.7886
$ws.Cells.Item(105, 9).Value = 2185.5642
$ws.Cells.Item(105, 11).Value = 2185.5642
$ws.Cells.Item(105, 13).Value = -438.5641999999998
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(55, 8).Value = 8111
$ws.Cells.Item(55, 9).Value = 8111
$ws.Cells.Item(55, 11).Value = 8111
$ws.Cells.Item(55, 13).Value = -7796
$ws.Cells.Item(97, 8).Value = 75000
$ws.Cells.Item(97, 9).Value = 75000
$ws.Cells.Item(97, 11).Value = 75000
$ws.Cells.Item(97, 13).Value = -74009
$ws.Cells.Item(122, 8).Value = 6632.8057
$ws.Cells.Item(122, 9).Value = 1626.3
$ws.Cells.Item(122, 11).Value = 4878.9
$ws.Cells.Item(122, 13).Value = -2428.9
$ws.Cells.Item(134, 8).Value = 3238.6086
$ws.Cells.Item(134, 9).Value = 2609.8948
$ws.Cells.Item(134, 11).Value = 7829.6844
$ws.Cells.Item(134, 13).Value = -5294.6844
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2920.8235
$ws.Cells.Item(5, 9).Value = 1763.4286
$ws.Cells.Item(5, 10).Value = 3731
$ws.Cells.Item(5, 11).Value = 5290.2858
$ws.Cells.Item(5, 12).Value = 11193
$ws.Cells.Item(5, 13).Value = -5178.2858
$ws.Cells.Item(5, 14).Value = -11417
$ws.Cells.Item(55, 8).Value = 46672710
$ws.Cells.Item(55, 10).Value = 6799.375
$ws.Cells.Item(55, 12).Value = 20398.125
$ws.Cells.Item(55, 14).Value = -20752.125
$ws.Cells.Item(64, 8).Value = 9518.5
$ws.Cells.Item(64, 9).Value = 5359.8
$ws.Cells.Item(64, 11).Value = 16079.4
$ws.Cells.Item(64, 13).Value = -15809.4
$ws.Cells.Item(67, 8).Value = 9518.5
$ws.Cells.Item(67, 9).Value = 5359.8
$ws.Cells.Item(67, 11).Value = 16079.4
$ws.Cells.Item(67, 13).Value = -15143.4
$ws.Cells.Item(107, 8).Value = 3760.389
$ws.Cells.Item(107, 9).Value = 840.6
$ws.Cells.Item(107, 10).Value = 4883.385
$ws.Cells.Item(107, 11).Value = 2521.8
$ws.Cells.Item(107, 12).Value = 14650.155
$ws.Cells.Item(107, 13).Value = -601.8000000000002
$ws.Cells.Item(107, 14).Value = -18490.155
$ws.Cells.Item(114, 8).Value = 10578.571
$ws.Cells.Item(114, 9).Value = 1583.3334
$ws.Cells.Item(114, 11).Value = 4750.0002
$ws.Cells.Item(114, 13).Value = -1496.0002
$ws.Cells.Item(117, 8).Value = 2872.2
$ws.Cells.Item(117, 10).Value = 2999.5
$ws.Cells.Item(117, 12).Value = 8998.5
$ws.Cells.Item(117, 14).Value = -15882.5
$ws.Cells.Item(126, 8).Value = 11343.333
$ws.Cells.Item(126, 9).Value = 2030
$ws.Cells.Item(126, 10).Value = 16000
$ws.Cells.Item(126, 11).Value = 6090
$ws.Cells.Item(126, 12).Value = 48000
$ws.Cells.Item(126, 13).Value = -1150
$ws.Cells.Item(126, 14).Value = -57880
$ws.Cells.Item(127, 8).Value = 11285
$ws.Cells.Item(127, 10).Value = 11285
$ws.Cells.Item(127, 12).Value = 33855
$ws.Cells.Item(127, 14).Value = -43775
$ws.Cells.Item(129, 8).Value = 1002906.6
$ws.Cells.Item(129, 9).Value = 1252036.2
$ws.Cells.Item(129, 10).Value = 6388
$ws.Cells.Item(129, 11).Value = 3756108.6
$ws.Cells.Item(129, 12).Value = 19164
$ws.Cells.Item(129, 13).Value = -3751108.6
$ws.Cells.Item(129, 14).Value = -29164
$ws.Cells.Item(132, 8).Value = 967.15
$ws.Cells.Item(132, 10).Value = 1213.1818
$ws.Cells.Item(132, 12).Value = 10918.6362
$ws.Cells.Item(132, 14).Value = -15978.6362
$ws.Cells.Item(135, 8).Value = 2920.8235
$ws.Cells.Item(135, 9).Value = 1763.4286
$ws.Cells.Item(135, 10).Value = 3731
$ws.Cells.Item(135, 11).Value = 15870.8574
$ws.Cells.Item(135, 12).Value = 33579
$ws.Cells.Item(135, 13).Value = -13335.8574
$ws.Cells.Item(135, 14).Value = -38649
$ws.Cells.Item(138, 8).Value = 25396.26
$ws.Cells.Item(138, 9).Value = 41999.89
$ws.Cells.Item(138, 10).Value = 17094.445
$ws.Cells.Item(138, 11).Value = 125999.67
$ws.Cells.Item(138, 12).Value = 51283.335
$ws.Cells.Item(138, 13).Value = -120859.67
$ws.Cells.Item(138, 14).Value = -61563.335
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2258.4412
$ws.Cells.Item(22, 10).Value = 2258.4412
$ws.Cells.Item(22, 12).Value = 2258.4412
$ws.Cells.Item(22, 14).Value = -2848.4412
$ws.Cells.Item(27, 8).Value = 2258.4412
$ws.Cells.Item(27, 10).Value = 2258.4412
$ws.Cells.Item(27, 12).Value = 2258.4412
$ws.Cells.Item(27, 14).Value = -2472.4412
$ws.Cells.Item(40, 8).Value = 4823.3335
$ws.Cells.Item(40, 9).Value = 5001.875
$ws.Cells.Item(40, 11).Value = 5001.875
$ws.Cells.Item(40, 13).Value = -4865.875
$ws.Cells.Item(138, 8).Value = 78929
$ws.Cells.Item(138, 10).Value = 78929
$ws.Cells.Item(138, 12).Value = 78929
$ws.Cells.Item(138, 14).Value = -89209
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(34, 8).Value = 4000
$ws.Cells.Item(34, 9).Value = 4000
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 4000
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = -3797
$ws.Cells.Item(34, 14).ClearContents()
$ws.Cells.Item(96, 8).Value = 2766.6667
$ws.Cells.Item(96, 9).Value = 1900
$ws.Cells.Item(96, 10).Value = 4500
$ws.Cells.Item(96, 11).Value = 1900
$ws.Cells.Item(96, 12).Value = 4500
$ws.Cells.Item(96, 13).Value = -527
$ws.Cells.Item(96, 14).Value = -7246
$ws.Cells.Item(136, 8).Value = 3086.4443
$ws.Cells.Item(136, 9).Value = 2847.25
$ws.Cells.Item(136, 11).Value = 8541.75
$ws.Cells.Item(136, 13).Value = -5991.75
